$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "26.952.86"
$ws.Range("E2").Value = "  -0.54%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.673.69"
$ws.Range("E3").Value = "  -0.26%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.08%  "

# Row 5 - BNB
Set-TextValue "D5" "214.71"
$ws.Range("E5").Value = "  -0.89%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  +1.69%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.05%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.250"
$ws.Range("E8").Value = "  -0.91%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -0.14%  "

# Row 10 - Solana
Set-TextValue "D10" "20.37"
$ws.Range("E10").Value = "  +0.64%  "

# Row 11 - TRON
Set-TextValue "D11" "0.0887"
$ws.Range("E11").Value = "  +0.02%  "

# Row 12 - WrappedliquidstakedEther2.0
Set-TextValue "D12" "1.910.58"
$ws.Range("E12").Value = "  -0.26%  "

# Row 13 - WrappedEther
Set-TextValue "D13" "1.681.32"
$ws.Range("E13").Value = "  +0.04%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -0.49%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  +0.46%  "

# Row 16 - Litecoin
Set-TextValue "D16" "65.64"
$ws.Range("E16").Value = "  -0.61%  "

# Row 17 - WrappedBTC
Set-TextValue "D17" "26.963.37"
$ws.Range("E17").Value = "  -0.56%  "

# Row 18 - BitcoinCash
Set-TextValue "D18" "236.32"
$ws.Range("E18").Value = "  -1.13%  "

# Row 19 - Chainlink
Set-TextValue "D19" "8.05"
$ws.Range("E19").Value = "  +3.77%  "

# Row 20 - ShibaInu
Set-TextValue "D20" "0.0₃0732"
$ws.Range("E20").Value = "  -0.90%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.11%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -1.00%  "

# Row 23 - Avalanche
Set-TextValue "D23" "9.18"
$ws.Range("E23").Value = "  -1.39%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -2.43%  "

# Row 25 - Monero
$ws.Range("E25").Value = "  -0.25%  "

# Row 26 - Cosmos
Set-TextValue "D26" "7.22"
$ws.Range("E26").Value = "  +1.10%  "

# Row 27 - EthereumClassic
Set-TextValue "D27" "15.97"
$ws.Range("E27").Value = "  -0.51%  "

# Row 28 - Stellar
Set-TextValue "D28" "0.112"
$ws.Range("E28").Value = "  -1.48%  "

# Row 29 - BinanceUSD
$ws.Range("E29").Value = "  +0.09%  "

# Row 30 - Hedera
$ws.Range("E30").Value = "  -0.49%  "

# Row 31 - PancakeSwap
$ws.Range("E31").Value = "  -0.48%  "

# Row 32 - Filecoin
Set-TextValue "D32" "3.32"
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - Maker
Set-TextValue "D33" "1.479.77"
$ws.Range("E33").Value = "  -0.23%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  +0.45%  "

# Row 35 - LidoDAOToken
Set-TextValue "D35" "1.68"
$ws.Range("E35").Value = "  +3.70%  "

# Row 36 - HuobiToken
Set-TextValue "D36" "2.41"

# Row 37 - ImmutableX
Set-TextValue "D37" "0.585"
$ws.Range("E37").Value = "  +1.28%  "

# Rows 38/39 - ARBITRUM and VeChain swap places
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D38" "0.0171"
$ws.Range("E38").Value = "  +1.13%  "

$ws.Range("B39").Value = "ARBITRUM"
$ws.Range("C39").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D39" "0.895"
$ws.Range("E39").Value = "  -0.95%  "

# Row 40 - FraxShare
$ws.Range("E40").Value = "  -2.99%  "

# Row 41 - WEMIXToken
$ws.Range("E41").Value = "  +4.97%  "

# Row 42 - PaxDollar
$ws.Range("E42").Value = "  +0.13%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +1.14%  "

# Row 44 - Aave
Set-TextValue "D44" "66.86"
$ws.Range("E44").Value = "  -0.06%  "

# Row 45 - RocketPoolETH
Set-TextValue "D45" "1.816.93"
$ws.Range("E45").Value = "  -0.39%  "

# Row 46 - TrustWalletToken
Set-TextValue "D46" "0.776"
$ws.Range("E46").Value = "  -0.32%  "

# Row 47 - Quant
Set-TextValue "D47" "90.43"
$ws.Range("E47").Value = "  +0.01%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.47%  "

# Row 49 - Algorand
$ws.Range("E49").Value = "  +0.90%  "

# Row 50 - Cronos
Set-TextValue "D50" "0.0508"
$ws.Range("E50").Value = "  -0.02%  "

# Row 51 - EnergySwap
Set-TextValue "D51" "7.69"
$ws.Range("E51").Value = "  -0.47%  "
